$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.37061877171383
$ws.Range("C2").Value = 9.884314789811398
$ws.Range("D2").Value = 14.39115976362049
$ws.Range("E2").Value = 15.39957656373883
$ws.Range("G2").Value = 3.661767390434881
$ws.Range("I2").Value = 22.97566110531809
$ws.Range("J2").Value = 9.061049849680234
$ws.Range("M2").Value = 18.67517686617916
$ws.Range("N2").Value = 18.20806076758938
$ws.Range("O2").Value = 26.4942358580312
$ws.Range("B3").Value = 13.9295752237302
$ws.Range("C3").Value = 9.501624798638456
$ws.Range("D3").Value = 14.37751138725411
$ws.Range("E3").Value = 15.41717945590063
$ws.Range("G3").Value = 3.664234663915335
$ws.Range("I3").Value = 23.06535142602514
$ws.Range("J3").Value = 9.081696179924258
$ws.Range("M3").Value = 18.54634220505655
$ws.Range("N3").Value = 18.26955369139376
$ws.Range("O3").Value = 26.53305918157813
$ws.Range("B4").Value = 13.65391967627984
$ws.Range("C4").Value = 9.26016709825886
$ws.Range("D4").Value = 14.3721014056666
$ws.Range("E4").Value = 15.430789734748
$ws.Range("G4").Value = 3.665830180166266
$ws.Range("I4").Value = 23.12551386557413
$ws.Range("J4").Value = 9.095195866346218
$ws.Range("M4").Value = 18.47017247404965
$ws.Range("N4").Value = 18.30915398998956
$ws.Range("O4").Value = 26.56353488202147
$ws.Range("B5").Value = 13.54054431278718
$ws.Range("C5").Value = 9.160292797815215
$ws.Range("D5").Value = 14.37064563839026
$ws.Range("E5").Value = 15.43704018746599
$ws.Range("G5").Value = 3.666500700716882
$ws.Range("I5").Value = 23.1513076832832
$ws.Range("J5").Value = 9.100904383524714
$ws.Range("M5").Value = 18.43989716337929
$ws.Range("N5").Value = 18.32575631872854
$ws.Range("O5").Value = 26.57761879196974
$ws.Range("B6").Value = 13.52166080978282
$ws.Range("C6").Value = 9.143624232046955
$ws.Range("D6").Value = 14.37044918549208
$ws.Range("E6").Value = 15.43812058316377
$ws.Range("G6").Value = 3.666613270284888
$ws.Range("I6").Value = 23.15566776461139
$ws.Range("J6").Value = 9.101864809692225
$ws.Range("M6").Value = 18.4349168898353
$ws.Range("N6").Value = 18.32854124051897
$ws.Range("O6").Value = 26.58005782531993
$ws.Range("B7").Value = 13.65239463863622
$ws.Range("C7").Value = 9.258825930907635
$ws.Range("D7").Value = 14.3720787384061
$ws.Range("E7").Value = 15.43087118016468
$ws.Range("G7").Value = 3.665839140617236
$ws.Range("I7").Value = 23.12585656342837
$ws.Range("J7").Value = 9.095272013588678
$ws.Range("M7").Value = 18.46976104139205
$ws.Range("N7").Value = 18.30937601071722
$ws.Range("O7").Value = 26.56371808807159
$ws.Range("B8").Value = 14.21965636437547
$ws.Range("C8").Value = 9.753795429394579
$ws.Range("D8").Value = 14.38583866856195
$ws.Range("E8").Value = 15.40506426410891
$ws.Range("G8").Value = 3.662601418183669
$ws.Range("I8").Value = 23.00552743593456
$ws.Range("J8").Value = 9.067998179697243
$ws.Range("M8").Value = 18.63016164736296
$ws.Range("N8").Value = 18.22888174426399
$ws.Range("O8").Value = 26.50624197091135
$ws.Range("B9").Value = 15.28643348655992
$ws.Range("C9").Value = 10.66710558272861
$ws.Range("D9").Value = 14.43628236233661
$ws.Range("E9").Value = 15.37670382622635
$ws.Range("G9").Value = 3.656888722264705
$ws.Range("I9").Value = 22.81011548319137
$ws.Range("J9").Value = 9.021024620991758
$ws.Range("M9").Value = 18.96678411976186
$ws.Range("N9").Value = 18.08559724661002
$ws.Range("O9").Value = 26.44637140841317
$ws.Range("B10").Value = 16.03378959100406
$ws.Range("C10").Value = 11.29637966987995
$ws.Range("D10").Value = 14.48747353121877
$ws.Range("E10").Value = 15.36944029447416
$ws.Range("G10").Value = 3.653075332850133
$ws.Range("I10").Value = 22.69147047215374
$ws.Range("J10").Value = 8.990457391337619
$ws.Range("M10").Value = 19.2258674305911
$ws.Range("N10").Value = 17.98911764158692
$ws.Range("O10").Value = 26.43479208348415
$ws.Range("B11").Value = 16.36432399703952
$ws.Range("C11").Value = 11.57243716054887
$ws.Range("D11").Value = 14.51378288162639
$ws.Range("E11").Value = 15.369081508583
$ws.Range("G11").Value = 3.651422935401769
$ws.Range("I11").Value = 22.64295054407105
$ws.Range("J11").Value = 8.977402857247903
$ws.Range("M11").Value = 19.34590921417255
$ws.Range("N11").Value = 17.94711707999212
$ws.Range("O11").Value = 26.43658872570796
$ws.Range("B12").Value = 16.48801888808896
$ws.Range("C12").Value = 11.67542638521036
$ws.Range("D12").Value = 14.5241752960937
$ws.Range("E12").Value = 15.36936865370393
$ws.Range("G12").Value = 3.6508089855047
$ws.Range("I12").Value = 22.62536484664629
$ws.Range("J12").Value = 8.972581382506499
$ws.Range("M12").Value = 19.39164730105387
$ws.Range("N12").Value = 17.93148276475426
$ws.Range("O12").Value = 26.43828581799944
$ws.Range("B13").Value = 16.46144609861593
$ws.Range("C13").Value = 11.65331581705322
$ws.Range("D13").Value = 14.52191807983703
$ws.Range("E13").Value = 15.36928800930674
$ws.Range("G13").Value = 3.65094068779293
$ws.Range("I13").Value = 22.62911714739733
$ws.Range("J13").Value = 8.973614352820544
$ws.Range("M13").Value = 19.3817848249006
$ws.Range("N13").Value = 17.93483788733565
$ws.Range("O13").Value = 26.43787509234943
$ws.Range("B14").Value = 16.37453054552949
$ws.Range("C14").Value = 11.58094158810725
$ws.Range("D14").Value = 14.51462929147696
$ws.Range("E14").Value = 15.36909665830631
$ws.Range("G14").Value = 3.65137218969685
$ws.Range("I14").Value = 22.64148794945801
$ws.Range("J14").Value = 8.977003748541811
$ws.Range("M14").Value = 19.34966667599981
$ws.Range("N14").Value = 17.94582542321091
$ws.Range("O14").Value = 26.4367079669083
$ws.Range("B15").Value = 16.32109748483239
$ws.Range("C15").Value = 11.5364065151374
$ws.Range("D15").Value = 14.51022049355962
$ws.Range("E15").Value = 15.36903451890183
$ws.Range("G15").Value = 3.651638028968888
$ws.Range("I15").Value = 22.64916812011272
$ws.Range("J15").Value = 8.979095728842502
$ws.Range("M15").Value = 19.33002897091555
$ws.Range("N15").Value = 17.95259078379079
$ws.Range("O15").Value = 26.43612549314865
$ws.Range("B16").Value = 16.01198841054382
$ws.Range("C16").Value = 11.27812621315521
$ws.Range("D16").Value = 14.48581455586067
$ws.Range("E16").Value = 15.36952297037525
$ws.Range("G16").Value = 3.653184972326574
$ws.Range("I16").Value = 22.69475139611843
$ws.Range("J16").Value = 8.991327624911625
$ws.Range("M16").Value = 19.21806353307589
$ws.Range("N16").Value = 17.99190037733958
$ws.Range("O16").Value = 26.43481691419651
$ws.Range("B17").Value = 15.81985705764225
$ws.Range("C17").Value = 11.11700553101057
$ws.Range("D17").Value = 14.47161300075647
$ws.Range("E17").Value = 15.37057669962633
$ws.Range("G17").Value = 3.654155014934493
$ws.Range("I17").Value = 22.72411437874284
$ws.Range("J17").Value = 8.999049133751754
$ws.Range("M17").Value = 19.14991217326959
$ws.Range("N17").Value = 18.01649838500654
$ws.Range("O17").Value = 26.43582434374252
$ws.Range("B18").Value = 15.70846549115759
$ws.Range("C18").Value = 11.02337693953752
$ws.Range("D18").Value = 14.46372943747174
$ws.Range("E18").Value = 15.37146006069236
$ws.Range("G18").Value = 3.654720710999203
$ws.Range("I18").Value = 22.7415161650794
$ws.Range("J18").Value = 9.003570428103391
$ws.Range("M18").Value = 19.11092084479796
$ws.Range("N18").Value = 18.03082433162827
$ws.Range("O18").Value = 26.43706869335424
$ws.Range("B19").Value = 15.67060227880354
$ws.Range("C19").Value = 10.9915143152459
$ws.Range("D19").Value = 14.46110925519243
$ws.Range("E19").Value = 15.37180678825361
$ws.Range("G19").Value = 3.654913579490348
$ws.Range("I19").Value = 22.74749610098593
$ws.Range("J19").Value = 9.005115025759897
$ws.Range("M19").Value = 19.0977556971894
$ws.Range("N19").Value = 18.03570542850572
$ws.Range("O19").Value = 26.43760416555537
$ws.Range("B20").Value = 15.84040198511867
$ws.Range("C20").Value = 11.13425668927492
$ws.Range("D20").Value = 14.47309533960071
$ws.Range("E20").Value = 15.37043583373087
$ws.Range("G20").Value = 3.654050950265251
$ws.Range("I20").Value = 22.72093552098321
$ws.Range("J20").Value = 8.998218879676376
$ws.Range("M20").Value = 19.15714576183357
$ws.Range("N20").Value = 18.01386148963145
$ws.Range("O20").Value = 26.43564827697508
$ws.Range("B21").Value = 16.40010052721133
$ws.Range("C21").Value = 11.60224224720328
$ws.Range("D21").Value = 14.51675856713292
$ws.Range("E21").Value = 15.36914138775359
$ws.Range("G21").Value = 3.651245127992412
$ws.Range("I21").Value = 22.6378329364597
$ws.Range("J21").Value = 8.976004892569087
$ws.Range("M21").Value = 19.35909319783014
$ws.Range("N21").Value = 17.94259078867696
$ws.Range("O21").Value = 26.43702318227075
$ws.Range("B22").Value = 16.75727557387345
$ws.Range("C22").Value = 11.89904650695414
$ws.Range("D22").Value = 14.54779659996787
$ws.Range("E22").Value = 15.37076061184269
$ws.Range("G22").Value = 3.649479978682298
$ws.Range("I22").Value = 22.58811386783067
$ws.Range("J22").Value = 8.962197693495868
$ws.Range("M22").Value = 19.49269906532248
$ws.Range("N22").Value = 17.89758674834097
$ws.Range("O22").Value = 26.44384799911583
$ws.Range("B23").Value = 16.56746774607485
$ws.Range("C23").Value = 11.74148852968376
$ws.Range("D23").Value = 14.5310038882887
$ws.Range("E23").Value = 15.36967106823383
$ws.Range("G23").Value = 3.650415813984213
$ws.Range("I23").Value = 22.61422837190499
$ws.Range("J23").Value = 8.96950191297308
$ws.Range("M23").Value = 19.42125375073984
$ws.Range("N23").Value = 17.92146247910033
$ws.Range("O23").Value = 26.43966310892684
$ws.Range("B24").Value = 15.83111651974843
$ws.Range("C24").Value = 11.12646054346979
$ws.Range("D24").Value = 14.47242429824807
$ws.Range("E24").Value = 15.37049865459532
$ws.Range("G24").Value = 3.654097972960448
$ws.Range("I24").Value = 22.72237106110664
$ws.Range("J24").Value = 8.998593981789853
$ws.Range("M24").Value = 19.15387486270865
$ws.Range("N24").Value = 18.01505305606097
$ws.Range("O24").Value = 26.43572580482661
$ws.Range("B25").Value = 15.0036650521918
$ws.Range("C25").Value = 10.42692513477017
$ws.Range("D25").Value = 14.42014011488211
$ws.Range("E25").Value = 15.38199251281198
$ws.Range("G25").Value = 3.65836646365585
$ws.Range("I25").Value = 22.85861712898386
$ws.Range("J25").Value = 9.033037802573537
$ws.Range("M25").Value = 18.87352647670668
$ws.Range("N25").Value = 18.08559724661002
$ws.Range("O25").Value = 26.45688820967
